$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---- Remove the two extra data rows (rows 3 and 4), keeping header + one data row ----
$ws1.Rows("3:4").Delete()

# ---- Update row 2 with the new contract data ----
$contractDate = Get-Date -Year 2021 -Month 4 -Day 2 -Hour 0 -Minute 0 -Second 0

$ws1.Range("A2").Value = "BCDID17"
$ws1.Range("B2").Value = "Army"
$ws1.Range("C2").Value = "BCDNO17"
$ws1.Range("D2").Value = $contractDate
$ws1.Range("E2").Value = "Hi, This is Description"
$ws1.Range("F2").Value = "this is category"
$ws1.Range("G2").Value = $contractDate
$ws1.Range("H2").Value = $contractDate
$ws1.Range("I2").Value = 5
$ws1.Range("J2").Value = $contractDate
$ws1.Range("K2").Value = "this is Incoterms"
$ws1.Range("L2").Value = 3
$ws1.Range("M2").Value = 5000
$ws1.Range("N2").Value = "Hi this is contents"
$ws1.Range("O2").Value = "Hi this is ny duties"
$ws1.Range("P2").Value = $contractDate
$ws1.Range("Q2").Value = "3 months"

# ---- Update the view selection state on Sheet1 ----
$ws1.Activate()
[void]$ws1.Range("E14").Select()
